$d = $word.ActiveDocument

# The first paragraph contains the placeholder ID text followed by a lone
# space run. Update the paragraph formatting (border + indent) and replace
# the two runs with a single run containing the new ID text.
$p1 = $d.Paragraphs.Item(1)

# Add a paragraph border (top/left/bottom/right, space = 5) to match the
# borders already used later in the document.
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

$p1.Format.LeftIndent = 11.25

$d.Content.Find.Execute("**ID__AFFARS_pgi_5349_topic_7__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_USAFA_PGI_5349_101__ID**", 2)
